# Corrigir pequenos erros de ortografia.

$d = $word.ActiveDocument

# 1. "...que ainda não estão " -> "...que ainda não está "
#    (concordância verbal: "um dos veterinários" é singular)
$d.Content.Find.Execute(
    "que ainda não estão ", $true, $false, $false, $false, $false,
    $true, 1, $false, "que ainda não está ", 2)

# 2. "5. Sistema válida CRMV." -> "5. Sistema valida CRMV."
#    (erro de ortografia: "valida", forma verbal, não leva acento)
$d.Content.Find.Execute(
    "5. Sistema válida CRMV.", $true, $false, $false, $false, $false,
    $true, 1, $false, "5. Sistema valida CRMV.", 2)
